$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for columns that get permuted across rows 2-24
$orig = @{}
$orig[2] = @{
    D = 44489
    K = 'Hayward'
    L = 'Primera'
    M = 300
    N = 26000
    O = 27000
    P = 26500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1472
    T = 18
}
$orig[3] = @{
    D = 45043
    K = 'Hayward'
    L = 'Segunda'
    M = 300
    N = 21000
    O = 22000
    P = 21500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1194
    T = 18
}
$orig[4] = @{
    D = 44629
    K = 'Hayward'
    L = 'Segunda'
    M = 300
    N = 17000
    O = 18000
    P = 17500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 972
    T = 18
}
$orig[5] = @{
    D = 44819
    K = 'Hayward'
    L = 'Primera'
    M = 300
    N = 17000
    O = 18000
    P = 17500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 1750
    T = 10
}
$orig[6] = @{
    D = 44991
    K = 'Hayward'
    L = 'Primera'
    M = 250
    N = 24000
    O = 25000
    P = 24500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1361
    T = 18
}
$orig[7] = @{
    D = 44323
    K = 'Hayward'
    L = 'Primera'
    M = 270
    N = 21000
    O = 22000
    P = 21500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1194
    T = 18
}
$orig[8] = @{
    D = 44418
    K = 'Hayward'
    L = 'Primera'
    M = 240
    N = 10000
    O = 11000
    P = 10500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 1050
    T = 10
}
$orig[9] = @{
    D = 45034
    K = 'Hayward'
    L = 'Primera'
    M = 250
    N = 25000
    O = 26000
    P = 25600
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1422
    T = 18
}
$orig[10] = @{
    D = 44307
    K = 'Hayward'
    L = 'Primera'
    M = 250
    N = 19000
    O = 20000
    P = 19500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1083
    T = 18
}
$orig[11] = @{
    D = 44602
    K = 'Hayward'
    L = 'Primera'
    M = 270
    N = 20000
    O = 21000
    P = 20500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1139
    T = 18
}
$orig[12] = @{
    D = 44263
    K = 'Hayward'
    L = 'Primera'
    M = 250
    N = 21000
    O = 22000
    P = 21500
    Q = '$/caja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1194
    T = 18
}
$orig[13] = @{
    D = 44784
    K = 'Hayward'
    L = 'Primera'
    M = 300
    N = 19000
    O = 20000
    P = 19500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1083
    T = 18
}
$orig[14] = @{
    D = 45002
    K = 'Hayward'
    L = 'Segunda'
    M = 300
    N = 24000
    O = 25000
    P = 24500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1361
    T = 18
}
$orig[15] = @{
    D = 45069
    K = 'Sin especificar'
    L = 'Primera'
    M = 370
    N = 19000
    O = 20000
    P = 19486
    Q = '$/bandeja 18 kilos'
    R = 'Región Metropolitana'
    S = 1083
    T = 18
}
$orig[16] = @{
    D = 44616
    K = 'Hayward'
    L = 'Segunda'
    M = 300
    N = 16000
    O = 17000
    P = 16500
    Q = '$/caja 18 kilos granel'
    R = 'Región de O''Higgins'
    S = 917
    T = 18
}
$orig[17] = @{
    D = 44614
    K = 'Hayward'
    L = 'Primera'
    M = 250
    N = 20000
    O = 21000
    P = 20500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1139
    T = 18
}
$orig[18] = @{
    D = 44789
    K = 'Hayward'
    L = 'Segunda'
    M = 250
    N = 19000
    O = 20000
    P = 19500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1083
    T = 18
}
$orig[19] = @{
    D = 44291
    K = 'Hayward'
    L = 'Primera'
    M = 200
    N = 17000
    O = 18000
    P = 17500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 972
    T = 18
}
$orig[20] = @{
    D = 44673
    K = 'Hayward'
    L = 'Especial'
    M = 400
    N = 14000
    O = 15000
    P = 14500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 1450
    T = 10
}
$orig[21] = @{
    D = 44487
    K = 'Hayward'
    L = 'Primera'
    M = 300
    N = 14000
    O = 15000
    P = 14500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 1450
    T = 10
}
$orig[22] = @{
    D = 44491
    K = 'Hayward'
    L = 'Primera'
    M = 300
    N = 14000
    O = 15000
    P = 14500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 1450
    T = 10
}
$orig[23] = @{
    D = 44656
    K = 'Hayward'
    L = 'Primera'
    M = 270
    N = 19000
    O = 20000
    P = 19500
    Q = '$/bandeja 18 kilos'
    R = 'Región de O''Higgins'
    S = 1083
    T = 18
}
$orig[24] = @{
    D = 44706
    K = 'Hayward'
    L = 'Primera'
    M = 400
    N = 9000
    O = 10000
    P = 9500
    Q = '$/bandeja 10 kilos'
    R = 'Región de O''Higgins'
    S = 950
    T = 10
}

# Apply permutation: row r gets values originally from row perm[r]
$perm = @{
    2 = 13
    3 = 19
    4 = 17
    5 = 12
    6 = 18
    7 = 7
    8 = 15
    9 = 24
    10 = 5
    11 = 2
    12 = 21
    13 = 10
    14 = 20
    15 = 16
    16 = 11
    17 = 22
    18 = 23
    19 = 6
    20 = 14
    21 = 8
    22 = 4
    23 = 3
    24 = 9
}

foreach ($r in $perm.Keys) {
    $src = $perm[$r]
    $vals = $orig[$src]
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 11).Value = $vals.K
    $ws.Cells.Item($r, 12).Value = $vals.L
    $ws.Cells.Item($r, 13).Value = $vals.M
    $ws.Cells.Item($r, 14).Value = $vals.N
    $ws.Cells.Item($r, 15).Value = $vals.O
    $ws.Cells.Item($r, 16).Value = $vals.P
    $ws.Cells.Item($r, 17).Value = $vals.Q
    $ws.Cells.Item($r, 18).Value = $vals.R
    $ws.Cells.Item($r, 19).Value = $vals.S
    $ws.Cells.Item($r, 20).Value = $vals.T
}